$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Row 7 corresponds to the 87c3a9a8-0836-4b55-851e-0bdffb843ea1 file handoff.
# This is a "Generate Report for Handoff" update: new handoff timestamps.

$wsOverview.Range("D7").Value = "2016-31-19 20:31:22"
$wsZhCn.Range("E7").Value = "2016-03-19 20:31:19"
$wsDeDe.Range("E7").Value = "2016-03-19 20:31:22"
